$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$ws.Range("F8").Value = "Wolf"
$ws.Range("D10").Value = "Stag"
$ws.Range("N10").Value = "Tufted Ground Squirrel"
$ws.Range("E12").Value = "Wolf"
$ws.Range("M12").Value = "Tufted Ground Squirrel"
$ws.Range("D14").Value = "Wolf"
$ws.Range("N14").Value = "Cobra Lily"
$ws.Range("G16").Value = "Wolf"
$ws.Range("D18").Value = "Indian Grey Mongoose"
$ws.Range("N18").Value = "Batfly"
$ws.Range("F24").Value = "Tiger"
$ws.Range("D26").Value = "Lucy "
$ws.Range("N26").Value = "Pitcher Plant"
$ws.Range("E28").Value = "Giant Squid "
$ws.Range("D30").Value = "Giant Squid "
$ws.Range("H32").Value = "Red-shanked Douc"
$ws.Range("I32").Value = "Red-shanked Douc"
$ws.Range("J32").Value = "Kob"
$ws.Range("D34").Value = "Rainbow Grasshopper"
$ws.Range("C35").Value = "Rainbow Grasshopper"
$ws.Range("E36").Value = "Nicobar Pigeon"
$ws.Range("D38").Value = "Nicobar Pigeon"
$ws.Range("F40").Value = "Nicobar Pigeon"
$ws.Range("D42").Value = "Halloween Crab"
$ws.Range("G48").Value = "Red-shanked Douc"
$ws.Range("K48").Value = "Great White Shark"
$ws.Range("D50").Value = "Mandarin Fish"
$ws.Range("E52").Value = "Wolf's Mona Monkey"
$ws.Range("M52").Value = "Great White Shark"
$ws.Range("F56").Value = "Red-shanked Douc"
$ws.Range("L56").Value = "Great White Shark"
$ws.Range("M60").Value = "Coral Snake"
